$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> D (price) and E (volume %) new text values.
# Only rows whose price actually changed get a D entry; every listed row gets an E entry.
$rows = @(
    @{ Row = 2;  D = "61.088.78";  E = "  +0.39%  " },
    @{ Row = 3;  D = "2.927.86";   E = "  +0.68%  " },
    @{ Row = 4;  E = "  -0.03%  " },
    @{ Row = 5;  D = "592.53";     E = "  +0.40%  " },
    @{ Row = 6;  D = "144.87";     E = "  -0.02%  " },
    @{ Row = 7;  E = "  +0.17%  " },
    @{ Row = 8;  D = "0.504";      E = "  -0.17%  " },
    @{ Row = 9;  E = "  +4.40%  " },
    @{ Row = 10; D = "0.142";      E = "  -0.78%  " },
    @{ Row = 11; E = "  -0.67%  " },
    @{ Row = 12; E = "  -0.44%  " },
    @{ Row = 13; D = "33.69";      E = "  +0.71%  " },
    @{ Row = 15; D = "3.413.30";   E = "  +0.68%  " },
    @{ Row = 16; D = "61.093.73";  E = "  +0.49%  " },
    @{ Row = 17; D = "6.72";       E = "  +0.18%  " },
    @{ Row = 18; D = "2.927.48";   E = "  +0.72%  " },
    @{ Row = 19; D = "433.43";     E = "  +0.96%  " },
    @{ Row = 20; E = "  -0.55%  " },
    @{ Row = 21; E = "  -0.46%  " },
    @{ Row = 22; E = "  +0.42%  " },
    @{ Row = 23; D = "81.45";      E = "  -0.59%  " },
    @{ Row = 24; D = "11.08";      E = "  +2.78%  " },
    @{ Row = 25; E = "  -0.44%  " },
    @{ Row = 26; E = "  -1.11%  " },
    @{ Row = 27; E = "  -0.02%  " },
    @{ Row = 28; E = "  -1.11%  " },
    @{ Row = 29; E = "  -0.37%  " },
    @{ Row = 30; E = "  -1.22%  " },
    @{ Row = 31; E = "  +2.30%  " },
    @{ Row = 32; D = "26.74";      E = "  +0.85%  " },
    @{ Row = 33; E = "  +0.05%  " },
    @{ Row = 34; D = "0.0₃0867";  E = "  +1.72%  " },
    @{ Row = 35; E = "  +0.32%  " },
    @{ Row = 36; D = "5.64";       E = "  +0.75%  " },
    @{ Row = 37; E = "  -1.03%  " },
    @{ Row = 38; E = "  -0.82%  " },
    @{ Row = 39; E = "  -0.50%  " },
    @{ Row = 40; D = "8.59";       E = "  -0.31%  " },
    @{ Row = 41; D = "41.98";      E = "  +4.77%  " },
    @{ Row = 42; D = "0.284";      E = "  -2.43%  " },
    @{ Row = 43; D = "373.03";     E = "  -0.03%  " },
    @{ Row = 44; E = "  -0.87%  " },
    @{ Row = 45; D = "2.710.08";   E = "  +0.39%  " },
    @{ Row = 46; D = "133.85";     E = "  +1.92%  " },
    @{ Row = 47; E = "  -0.04%  " },
    @{ Row = 48; D = "23.80";      E = "  -0.98%  " },
    @{ Row = 49; E = "  -1.19%  " },
    @{ Row = 50; E = "  -1.52%  " },
    @{ Row = 51; E = "  -0.56%  " }
)

# Price strings that look like a plain decimal number (single '.') would be
# auto-coerced to a numeric value by Excel when assigned via .Value (e.g.
# "592.53" -> 592.53). Force those specific cells to keep Text formatting so
# the literal string is preserved instead of becoming a float. Prices that
# contain more than one '.' (thousand separators) or other non-numeric
# characters are already safe and do not need this treatment.
$numericLookingRows = @(5, 6, 8, 10, 13, 17, 19, 23, 24, 32, 36, 40, 41, 42, 43, 46, 48)

foreach ($entry in $rows) {
    $r = $entry.Row
    if ($entry.ContainsKey("D")) {
        $cellD = $ws.Cells.Item($r, 4)
        if ($numericLookingRows -contains $r) {
            $cellD.NumberFormat = "@"
        }
        $cellD.Value = $entry.D
    }
    $ws.Cells.Item($r, 5).Value = $entry.E
}
